$d = $word.ActiveDocument

# --- 1) Merge the three split/proof-marked runs back into single runs ---
# These replacements keep the visible text identical; Word's COM Find/Replace
# (replacing text with itself) collapses the run that previously spanned
# <w:r>/<w:proofErr>/<w:r>/<w:proofErr>/<w:r> into one simple <w:r>.

$pairs = @(
    "Token zu sha String Umwandeln",
    "Die nötigen Informationen für die Entschlüsselung and den DB Queries bekommen",
    "Die Verschlüsselung an allen Queries s etablieren"
)

foreach ($text in $pairs) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

# --- 2) Append 13 blank paragraphs (the last one will then be turned into
#        the "Sendgrid" paragraph below), giving 12 blank paragraphs plus
#        one final paragraph with text split across 3 runs ---

for ($i = 0; $i -lt 13; $i++) {
    $end = $d.Content
    $end.Collapse(0)
    $end.Text = "`r"
}

$last = $d.Content
$last.Collapse(0)
$last.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>S</w:t></w:r><w:r><w:t>endgrid</w:t></w:r><w:r><w:t>: Mail API (Password via PM)</w:t></w:r></w:p>")

Write-Host "done"
